# Trade #12 closed at 2026-02-17 13:34:51 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: update aggregate stats after the new closed trade
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.79   # Current Capital
$summary.Range("B4").Value = -0.21     # Total P&L $
$summary.Range("B6").Value = 12        # Total Trades
$summary.Range("B8").Value = 9         # Losing Trades
$summary.Range("B9").Value = 25        # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.79000000000001  # Capital
$status.Range("D4").Value = 12                 # Trades
$status.Range("E4").Value = -0.21              # P&L $
$status.Range("F4").Value = -0.21              # P&L %
$status.Range("G4").Value = 25                 # Win Rate %

# ---------------------------------------------------------------------------
# Helper: write one new trade row across a trades-log style sheet
# (columns A..Q) while keeping the date column ("B") as plain text so
# Excel doesn't auto-convert the "yyyy-mm-dd" string into a date serial.
# ---------------------------------------------------------------------------
function Add-TradeRow {
    param($sheet, $rowIndex)

    $sheet.Cells.Item($rowIndex, 1).Value = 12
    $sheet.Cells.Item($rowIndex, 2).Value = "'2026-02-17"
    $sheet.Cells.Item($rowIndex, 3).Value = "13:34:45"
    $sheet.Cells.Item($rowIndex, 4).Value = "MarketMaking"
    $sheet.Cells.Item($rowIndex, 5).Value = "DOWN"
    $sheet.Cells.Item($rowIndex, 6).Value = 0.9
    $sheet.Cells.Item($rowIndex, 7).Value = 0.88
    $sheet.Cells.Item($rowIndex, 8).Value = "CLOSED"
    $sheet.Cells.Item($rowIndex, 9).Value = -2.2222
    $sheet.Cells.Item($rowIndex, 10).Value = -0.02
    $sheet.Cells.Item($rowIndex, 11).Value = 99.79000000000001
    $sheet.Cells.Item($rowIndex, 12).Value = 0
    $sheet.Cells.Item($rowIndex, 13).Value = 0
    $sheet.Cells.Item($rowIndex, 14).Value = 0.6
    $sheet.Cells.Item($rowIndex, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item($rowIndex, 16).Value = "early_exit"
    $sheet.Cells.Item($rowIndex, 17).Value = 0.13
}

# ---------------------------------------------------------------------------
# All Trades sheet: append new trade #12 as row 13
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 13

# ---------------------------------------------------------------------------
# MarketMaking sheet: append the same new trade #12 as row 13
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 13

Write-Host "Applied trade #12 close update."
